$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43
$ws.Cells.Item($row, 1).Value = 42
$ws.Cells.Item($row, 2).Value = 67
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 9
$ws.Cells.Item($row, 5).Value = 20
$ws.Cells.Item($row, 6).Value = 77
$ws.Cells.Item($row, 7).Value = 97
